$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.703.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.30%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.600.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.24%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  +0.29%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''211.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.06%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '''  -0.60%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.24%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.0619'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.25%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  +0.73%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''19.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +0.55%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.0845'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +0.97%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''1.825.12'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +0.25%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.614.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.52%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = '''  +0.50%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  +0.25%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''65.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +1.14%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''26.681.42'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +0.27%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''0.0₃0758'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +3.60%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = '''Chainlink'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = '''7.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +4.08%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = '''Dai'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = '''1.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.29%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''209.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.31%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +0.49%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''2.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -0.08%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''8.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +0.58%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''143.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -1.77%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +0.26%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''7.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +0.18%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D29").Value = '''15.35'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.62%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''0.0521'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +2.59%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = '''  -0.16%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +0.74%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = '''  +1.33%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.294.63'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +0.83%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -5.35%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  +1.12%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  +0.42%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -0.06%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.09'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +19.40%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.827'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.76%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -1.28%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -0.12%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -0.40%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''63.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -1.63%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.736.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +0.25%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''91.09'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +1.52%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''1.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -1.35%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = '''  -1.59%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  +0.55%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  +0.34%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''7.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -0.92%  '
$ws.Range("E51").Style = "Normal"
